# Update dSF (column F) values as part of a data repull / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F5").Value = -5
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -4
$ws.Range("F13").Value = -12
